$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = -7.676
    6   = -8.169
    7   = -7.842000000000001
    16  = -8.529
    20  = -7.672999999999999
    28  = -8.262
    29  = -7.292
    32  = -7.98
    40  = -7.883999999999999
    46  = -8.202000000000002
    51  = -8.457000000000001
    52  = -8.115
    57  = -8.06
    59  = -8.123999999999999
    62  = -7.946000000000001
    66  = -7.436
    73  = -7.789
    74  = -7.904999999999999
    92  = -6.536999999999999
    100 = -8.295
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
